$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the data table with the 2023 column (K) ---
# Copy the formatting of the current last data column (J) into the new
# column K for the header row and the three data rows, then fill in the
# 2023 figures.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("J4:J6").Copy()
$ws.Range("K4:K6").PasteSpecial(-4122)

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 3614
$ws.Range("K5").Value = 1301
$ws.Range("K6").Value = 2313

# Column K is now the right-hand edge of the table, so it needs the
# closing vertical border that used to not be necessary on column J.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2

# Keep the column widths consistent across the (now wider) data area.
$ws.Range("B1:O1").EntireColumn.ColumnWidth = 7.877604166666667

$excel.CutCopyMode = 0
